# Append several new paragraphs at the end of the document body, after the
# existing "This is the file to be tested" paragraph. Each new paragraph is
# inserted as a tiny WordprocessingML package via Range.InsertXML so that
# run boundaries (e.g. the two separate runs in "My fifth "/"change") are
# produced exactly rather than merged together.

$d = $word.ActiveDocument

function Append-ParagraphXml {
    param([string]$InnerParagraphXml)

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $InnerParagraphXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData>' +
           '</pkg:part>' +
           '</pkg:package>'

    $endRange = $d.Content
    $endRange.Collapse(0)   # wdCollapseEnd
    $endRange.InsertXML($pkg)
}

Append-ParagraphXml '<w:p><w:r><w:t>Now is the time for jillet blue 2</w:t></w:r></w:p>'
Append-ParagraphXml '<w:p><w:r><w:t>My third change</w:t></w:r></w:p>'
Append-ParagraphXml '<w:p><w:r><w:t>My fourth change</w:t></w:r></w:p>'
Append-ParagraphXml '<w:p><w:r><w:t>Also my fourth change</w:t></w:r></w:p>'
Append-ParagraphXml '<w:p><w:r><w:t xml:space="preserve">My fifth </w:t></w:r><w:r><w:t>change</w:t></w:r></w:p>'
Append-ParagraphXml '<w:p><w:r><w:t>My six change     also my six but same line</w:t></w:r></w:p>'
Append-ParagraphXml '<w:p><w:r><w:t>My six on next line</w:t></w:r></w:p>'

Write-Output $d.Paragraphs.Count
